$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 411, shifting existing rows 411-471 down to 412-472.
$ws.Rows(411).Insert()

# Populate the newly inserted row 411 with the new data record.
$ws.Range("A411").Value = 4
$ws.Range("B411").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C411").Value = "Los Lagos"
$ws.Range("D411").Value = 45127
$ws.Range("E411").Value = 10
$ws.Range("F411").Value = 100112043
$ws.Range("G411").Value = "Pepino ensalada"
$ws.Range("H411").Value = "Sin especificar"
$ws.Range("I411").Value = "Primera"
$ws.Range("J411").Value = 120
$ws.Range("K411").Value = 17000
$ws.Range("L411").Value = 17000
$ws.Range("M411").Value = 17000
$ws.Range("N411").Value = "`$/caja 60 unidades"
$ws.Range("O411").Value = "Región de Arica y Parinacota"
$ws.Range("P411").Value = 283
$ws.Range("Q411").Value = 60
$ws.Range("R411").Value = "Hortaliza"
